# Apply cryptos list update (prices & 1h volume % changes)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.717.83"
$ws.Range("E2").Value = "  +3.26%  "
$ws.Range("D3").Value = "3.695.12"
$ws.Range("E3").Value = "  +6.91%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.04"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.83"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").Value = "3.687.63"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.616"
$ws.Range("E8").Value = "  +3.59%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("E11").Value = "  +4.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.94"
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000287"
$ws.Range("E13").Value = "  +1.75%  "
$ws.Range("D14").Value = "4.296.41"
$ws.Range("E14").Value = "  +7.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "680.79"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "9.04"
$ws.Range("E16").Value = "  +3.82%  "
$ws.Range("D17").Value = "3.698.46"
$ws.Range("E17").Value = "  +7.25%  "
$ws.Range("D18").Value = "71.832.09"
$ws.Range("E18").Value = "  +3.34%  "
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.10"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("E21").Value = "  +2.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.37"
$ws.Range("E22").Value = "  +18.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.946"
$ws.Range("E23").Value = "  +3.81%  "
$ws.Range("E24").Value = "  +4.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "103.39"
$ws.Range("E25").Value = "  +2.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.04"
$ws.Range("E26").Value = "  +3.01%  "
$ws.Range("E27").Value = "  +4.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.33"
$ws.Range("E28").Value = "  +5.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.72"
$ws.Range("E29").Value = "  +5.45%  "
$ws.Range("E30").Value = "  +5.51%  "
$ws.Range("E31").Value = "  +6.23%  "
$ws.Range("E32").Value = "  +10.97%  "
$ws.Range("E33").Value = "  +1.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "567.14"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  +3.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.45"
$ws.Range("E36").Value = "  +2.31%  "
$ws.Range("D37").Value = "3.746.61"
$ws.Range("E37").Value = "  +2.80%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("E39").Value = "  +2.67%  "
$ws.Range("D40").Value = "0.0₃0776"
$ws.Range("E40").Value = "  +3.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.67"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.47"
$ws.Range("E42").Value = "  +4.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.81"
$ws.Range("E43").Value = "  +3.85%  "
$ws.Range("E44").Value = "  +9.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.352"
$ws.Range("E45").Value = "  +4.36%  "
$ws.Range("E46").Value = "  +8.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.39"
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("E48").Value = "  +3.01%  "
$ws.Range("E49").Value = "  +1.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.04"
$ws.Range("E51").Value = "  +3.55%  "
